$wb = $excel.ActiveWorkbook

# Remember the workbook's currently active sheet so we can restore selection at the end
$origActiveSheet = $excel.ActiveSheet

# The new daily ranking sheet goes right after the most recent existing date sheet ("2025-09-06")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "2025-09-07"

# Build the full rank/title/author/latest_episode table (header + 50 ranked rows) as a 2D array
$arr = New-Object 'object[,]' 51,4
$arr[0,0] = 'rank'
$arr[0,1] = 'title'
$arr[0,2] = 'author'
$arr[0,3] = 'latest_episode'
$arr[1,0] = 1
$arr[1,1] = 'クセ強彼女は床にいざなう'
$arr[1,2] = '須河篤志(著者)'
$arr[1,3] = '第15話前半'
$arr[2,0] = 2
$arr[2,1] = 'リビルドワールド'
$arr[2,2] = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$arr[2,3] = '第72話④'
$arr[3,0] = 3
$arr[3,1] = 'ぽんドロイド！ はまさん'
$arr[3,2] = 'はれやまはれぞう(著者)'
$arr[3,3] = '第8話'
$arr[4,0] = 4
$arr[4,1] = '氷結令嬢さまをフォローしたら、メチャメチャ溺愛されてしまった件@comic'
$arr[4,2] = '漫画：ハレノチアメ 原作：愛坂タカト キャラクター原案：Bcoca'
$arr[4,3] = 'アリシア様セクシーショット集（担当編集厳選）'
$arr[5,0] = 5
$arr[5,1] = 'まったく最近の探偵ときたら'
$arr[5,2] = '五十嵐正邦(著者)'
$arr[5,3] = '第115話'
$arr[6,0] = 6
$arr[6,1] = 'ダメ人間の愛しかた'
$arr[6,2] = '岩葉(著者)'
$arr[6,3] = '第19話後編　ダメ人間と新生活の彼女'
$arr[7,0] = 7
$arr[7,1] = '田舎の黒ギャルJKと結婚しました'
$arr[7,2] = 'カヅチ(著者)'
$arr[7,3] = '第19話前半'
$arr[8,0] = 8
$arr[8,1] = '理想の彼女'
$arr[8,2] = 'もりまりも(著者)'
$arr[8,3] = '第25話'
$arr[9,0] = 9
$arr[9,1] = '女友達は頼めば意外とヤらせてくれる'
$arr[9,2] = 'ろくろ(漫画) 鏡遊(原作)'
$arr[9,3] = '第24話'
$arr[10,0] = 10
$arr[10,1] = '無敵商人の異世界成り上がり物語 ～現代の製品を自在に取り寄せるスキルがあるので異世界では楽勝です～'
$arr[10,2] = '隆原ヒロタ(漫画) 青山有(原作) ぷきゅのすけ(キャラクターデザイン)'
$arr[10,3] = '第36話'
$arr[11,0] = 11
$arr[11,1] = 'ミルク搾りハンターの異世界搾乳記～農家の冴えない男があらゆる種族の地区Bを弄び虜にする～'
$arr[11,2] = '空詠大智(著者)'
$arr[11,3] = '第18話後半'
$arr[12,0] = 12
$arr[12,1] = '王子様の友達'
$arr[12,2] = 'すけろく(著者)'
$arr[12,3] = '第30話'
$arr[13,0] = 13
$arr[13,1] = 'リアリスト魔王による聖域なき異世界改革'
$arr[13,2] = '鈴木マナツ(漫画) 羽田遼亮(原作) ゆーげん(キャラクターデザイン) ひたきゆう(キャラクターデザイン)'
$arr[13,3] = '第68幕④'
$arr[14,0] = 14
$arr[14,1] = '転生してあらゆるモノに好かれながら異世界で好きな事をして生きて行く'
$arr[14,2] = '都尾琉(漫画) 御峰。(原作)'
$arr[14,3] = '第28話①'
$arr[15,0] = 15
$arr[15,1] = '王立魔術学院の鬼畜講師'
$arr[15,2] = '実々みみず(漫画) 急川回レ(原作) zunta(キャラクターデザイン)'
$arr[15,3] = '第20話'
$arr[16,0] = 16
$arr[16,1] = 'まんきつしたい常連さん'
$arr[16,2] = 'しんみりん(著者)'
$arr[16,3] = '第47話後編'
$arr[17,0] = 17
$arr[17,1] = '魔のものたちは企てる'
$arr[17,2] = '加藤拓弐(原作) ガしガし(作画)'
$arr[17,3] = '第29話'
$arr[18,0] = 18
$arr[18,1] = 'ダークサモナーとデキている'
$arr[18,2] = '車王(著者)'
$arr[18,3] = '第75話'
$arr[19,0] = 19
$arr[19,1] = '愚かな天使は悪魔と踊る'
$arr[19,2] = 'アズマサワヨシ(著者)'
$arr[19,3] = '第101話④'
$arr[20,0] = 20
$arr[20,1] = '今日から僕は、彼女の✕✕を解消する。'
$arr[20,2] = 'コアヤアコ(著者)'
$arr[20,3] = '第3話前半'
$arr[21,0] = 21
$arr[21,1] = 'いとこのこ'
$arr[21,2] = 'いぬちく(著者)'
$arr[21,3] = '単行本第4巻発売情報!!＆コメント企画発表!!'
$arr[22,0] = 22
$arr[22,1] = 'ワンパンマン'
$arr[22,2] = '原作/ＯＮＥ 作画/村田雄介'
$arr[22,3] = '210撃目'
$arr[23,0] = 23
$arr[23,1] = '最強ボイン・ランキング～ド陰キャニートのムスコ無双～'
$arr[23,2] = '衝撃の平山(著者)'
$arr[23,3] = '読切'
$arr[24,0] = 24
$arr[24,1] = 'このヒーラー、めんどくさい'
$arr[24,2] = '丹念に発酵(著者)'
$arr[24,3] = '特別編：ポーションの中身'
$arr[25,0] = 25
$arr[25,1] = '豚のレバーは加熱しろ'
$arr[25,2] = 'みなみ(漫画) 逆井卓馬(原作) 遠坂あさぎ(キャラクターデザイン)'
$arr[25,3] = '第43話②'
$arr[26,0] = 26
$arr[26,1] = 'きみの願いが叶うまで'
$arr[26,2] = '浅月のりと(著者)'
$arr[26,3] = '第5話-1'
$arr[27,0] = 27
$arr[27,1] = 'よくわからないけれど異世界に転生していたようです'
$arr[27,2] = '内々けやき あし カオミン'
$arr[27,3] = '第138話 よくわからないけれど逃亡中みたいです（２）'
$arr[28,0] = 28
$arr[28,1] = 'ずっと好きだった幼馴染と付き合い始めたら一途ビッチの性欲ジャンキーだったんだがどうすりゃいいですか？'
$arr[28,2] = '原作：トラ子猫 作画：あらいぐま'
$arr[28,3] = '第4話'
$arr[29,0] = 29
$arr[29,1] = '異世界のんびり農家'
$arr[29,2] = '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)'
$arr[29,3] = 'コミックス告知イラスト'
$arr[30,0] = 30
$arr[30,1] = 'ステータス・オール∞（インフィニティ） ∞使いの最強能力者、異世界を自由気ままに暮らします！'
$arr[30,2] = '漫画：三津屋みやこ 原作：八又ナガト'
$arr[30,3] = '第66話'
$arr[31,0] = 31
$arr[31,1] = '異種族追放コンカフェ'
$arr[31,2] = '佐々木マサヒト(著者)'
$arr[31,3] = '第18話-2'
$arr[32,0] = 32
$arr[32,1] = '小林さんちのメイドラゴン'
$arr[32,2] = 'クール教信者'
$arr[32,3] = '第149話'
$arr[33,0] = 33
$arr[33,1] = '器用貧乏、城を建てる～開拓学園の劣等生なのに、上級職のスキルと魔法がすべて使えます～＠COMIC'
$arr[33,2] = '原作：佐藤謙羊 漫画：スガン'
$arr[33,3] = '第23話②「自爆スイッチは押されたい」'
$arr[34,0] = 34
$arr[34,1] = '十歳の最強魔導師'
$arr[34,2] = '猫月 天乃聖樹'
$arr[34,3] = '第4話（前編）'
$arr[35,0] = 35
$arr[35,1] = '理想のヒモ生活'
$arr[35,2] = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$arr[35,3] = '第87話　その1'
$arr[36,0] = 36
$arr[36,1] = '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。'
$arr[36,2] = '六志麻あさ 業務用餅 kisui'
$arr[36,3] = '第７１話'
$arr[37,0] = 37
$arr[37,1] = 'おっさん、転生して天才役者になる'
$arr[37,2] = '芽々ノ圭(漫画) ほえ太郎(原作)'
$arr[37,3] = '第38話④'
$arr[38,0] = 38
$arr[38,1] = 'アザミヤコを好きになる'
$arr[38,2] = 'ユニティコング(原作) ツノニガウ(作画)'
$arr[38,3] = '描き下ろしお題大募集！！【第２巻発売記念コメント企画】'
$arr[39,0] = 39
$arr[39,1] = '悪役令嬢、俺'
$arr[39,2] = '弥南せいら(漫画) 猫屋敷のあ(原作)'
$arr[39,3] = '第2話-2'
$arr[40,0] = 40
$arr[40,1] = '王様ランキング200話～'
$arr[40,2] = '十日草輔（とおかそうすけ）'
$arr[40,3] = '第266話'
$arr[41,0] = 41
$arr[41,1] = '経験値貯蓄でのんびり傷心旅行 ～勇者と恋人に追放された戦士の無自覚ざまぁ～'
$arr[41,2] = '奏ヨシキ(著者) 徳川レモン(原作) riritto(キャラクターデザイン)'
$arr[41,3] = '第38話-2'
$arr[42,0] = 42
$arr[42,1] = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$arr[42,2] = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$arr[42,3] = '第82話その2'
$arr[43,0] = 43
$arr[43,1] = '陰キャの俺が席替えでS級美少女に囲まれたら秘密の関係が始まった。'
$arr[43,2] = '星野 星野(原作) バラマツヒトミ(漫画) 黒兎 ゆう(キャラクター原案)'
$arr[43,3] = '第5話'
$arr[44,0] = 44
$arr[44,1] = '無能の中の無能王子　スキル【無能】を授かりましたが、周りの女性は【傾国】【傾城】【奸婦】【毒婦】【悪婦】【妖婦】とかです'
$arr[44,2] = '漫画/一夢 原作/福朗 キャラクター原案/菊池政治'
$arr[44,3] = 'chapter12【7話①】'
$arr[45,0] = 45
$arr[45,1] = '修羅幼女の英雄譚～半端者と言われた傭兵、幼女に転生して成り上がる～'
$arr[45,2] = '作画：むらたん 原作：沙城流'
$arr[45,3] = '第9話(3)'
$arr[46,0] = 46
$arr[46,1] = '怪異部～M県Y市の怪現象について～'
$arr[46,2] = 'さりい・Ｂ(著者)'
$arr[46,3] = 'File.9'
$arr[47,0] = 47
$arr[47,1] = 'ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～'
$arr[47,2] = '神原絵理華(漫画) 一森一輝(原作)'
$arr[47,3] = '第19話②'
$arr[48,0] = 48
$arr[48,1] = '終末ツーリング'
$arr[48,2] = 'さいとー栄(著者)'
$arr[48,3] = '第50話 三沢基地　その６①'
$arr[49,0] = 49
$arr[49,1] = '陰キャの僕に罰ゲームで告白してきたはずのギャルが、どう見ても僕にベタ惚れです'
$arr[49,2] = '神奈なごみ(漫画) 結石(原作) かがちさく(キャラクター原案)'
$arr[49,3] = '第26.5話'
$arr[50,0] = 50
$arr[50,1] = '最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～'
$arr[50,2] = 'シオヤマ琴 鳥羽田 航 トモゼロ'
$arr[50,3] = '第79話 ポリリズム'

$newSheet.Range("A1:D51").Value = $arr

# Match the header-row formatting used on every other daily ranking sheet: bold, thin box border,
# centered horizontally and top-aligned vertically
$hdr = $newSheet.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# Restore original selection/active sheet (adding a sheet shouldn't change the user's active tab)
$origActiveSheet.Activate()
